$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-FormulaKeepFormat($range, $formula) {
    # Preserve the cell's existing number format / style (e.g. quote-prefix)
    # across a content replacement by round-tripping formats through a
    # scratch cell via copy/paste-special.
    $scratch = $ws.Range("ZZ1")
    $range.Copy()
    $scratch.PasteSpecial(-4122)   # xlPasteFormats
    $range.Formula = $formula
    $scratch.Copy()
    $range.PasteSpecial(-4122)     # xlPasteFormats
    $scratch.Clear()
}

# D16 held literal (quote-prefixed) text "=TRUE AND FALSE" -> now a real AND() formula
Set-FormulaKeepFormat $ws.Range("D16") "=AND(TRUE, FALSE)"

# D17 held literal (quote-prefixed) text "=TRUE OR FALSE" -> now a real OR() formula
Set-FormulaKeepFormat $ws.Range("D17") "=OR(TRUE, FALSE)"

# D20 formula changes from AND(TRUE,TRUE,FALSE) to XOR(TRUE, FALSE)
$ws.Range("D20").Formula = "=XOR(TRUE, FALSE)"
